$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4453.385
$ws.Range("J86").Value = 5928.5713
$ws.Range("L86").Value = 5928.5713
$ws.Range("N86").Value = -8174.5713
$ws.Range("H88").Value = 2607.7144
$ws.Range("I88").Value = 1040
$ws.Range("J88").Value = 3783.5
$ws.Range("K88").Value = 1040
$ws.Range("L88").Value = 3783.5
$ws.Range("M88").Value = -634
$ws.Range("N88").Value = -4595.5
$ws.Range("H89").Value = 4453.385
$ws.Range("J89").Value = 5928.5713
$ws.Range("L89").Value = 29642.8565
$ws.Range("N89").Value = -40874.85649999999
$ws.Range("H91").Value = 2607.7144
$ws.Range("I91").Value = 1040
$ws.Range("J91").Value = 3783.5
$ws.Range("K91").Value = 1040
$ws.Range("L91").Value = 3783.5
$ws.Range("M91").Value = 364
$ws.Range("N91").Value = -6591.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 25000
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H122").Value = 373605.72
$ws.Range("I122").Value = 502667.84
$ws.Range("K122").Value = 1508003.52
$ws.Range("M122").Value = -1505553.52
$ws.Range("H132").Value = 2103.9443
$ws.Range("I132").Value = 2117.0625
$ws.Range("K132").Value = 6351.1875
$ws.Range("M132").Value = -3821.1875

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1250
$ws.Range("I20").Value = 1250
$ws.Range("K20").Value = 1250
$ws.Range("M20").Value = -1003
$ws.Range("H105").Value = 3923.3044
$ws.Range("I105").Value = 2916.85
$ws.Range("J105").Value = 10633
$ws.Range("K105").Value = 2916.85
$ws.Range("L105").Value = 10633
$ws.Range("M105").Value = -1169.85
$ws.Range("N105").Value = -14127
$ws.Range("H132").Value = 106834.75
$ws.Range("J132").Value = 106834.75
$ws.Range("L132").Value = 106834.75
$ws.Range("N132").Value = -116954.75
$ws.Range("H134").Value = 1672.9375
$ws.Range("I134").Value = 1667.8
$ws.Range("K134").Value = 5003.4
$ws.Range("M134").Value = -2468.4

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 35000
$ws.Range("I59").Value = 35000
$ws.Range("K59").Value = 35000
$ws.Range("M59").Value = -33855
$ws.Range("H105").Value = 2394
$ws.Range("I105").Value = 1064.75
$ws.Range("J105").Value = 4166.3335
$ws.Range("K105").Value = 1064.75
$ws.Range("L105").Value = 4166.3335
$ws.Range("M105").Value = 682.25
$ws.Range("N105").Value = -7660.3335

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 12500057
$ws.Range("I7").Value = 50000030
$ws.Range("J7").Value = 65.5
$ws.Range("K7").Value = 150000090
$ws.Range("L7").Value = 196.5
$ws.Range("M7").Value = -149999978
$ws.Range("N7").Value = -420.5
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H131").Value = 1533.2325
$ws.Range("I131").Value = 1794.5
$ws.Range("J131").Value = 1520.4878
$ws.Range("K131").Value = 5383.5
$ws.Range("L131").Value = 4561.463400000001
$ws.Range("M131").Value = -343.5
$ws.Range("N131").Value = -14641.4634

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 725
$ws.Range("I43").Value = 725
$ws.Range("K43").Value = 725
$ws.Range("M43").Value = -574
$ws.Range("H80").Value = 4786.316
$ws.Range("I80").Value = 3614
$ws.Range("J80").Value = 6088.8887
$ws.Range("K80").Value = 3614
$ws.Range("L80").Value = 6088.8887
$ws.Range("M80").Value = -2616
$ws.Range("N80").Value = -8084.8887
$ws.Range("H83").Value = 4786.316
$ws.Range("I83").Value = 3614
$ws.Range("J83").Value = 6088.8887
$ws.Range("K83").Value = 18070
$ws.Range("L83").Value = 30444.4435
$ws.Range("M83").Value = -13078
$ws.Range("N83").Value = -40428.4435
$ws.Range("H132").Value = 2613.7693
$ws.Range("I132").Value = 907.7
$ws.Range("J132").Value = 8300.666999999999
$ws.Range("K132").Value = 2723.1
$ws.Range("L132").Value = 24902.001
$ws.Range("M132").Value = -193.1000000000004
$ws.Range("N132").Value = -29962.001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5656.143
$ws.Range("I16").Value = 2718.8
$ws.Range("J16").Value = 12999.5
$ws.Range("K16").Value = 2718.8
$ws.Range("L16").Value = 12999.5
$ws.Range("M16").Value = -2548.8
$ws.Range("N16").Value = -13339.5
$ws.Range("H61").Value = 5206.8184
$ws.Range("I61").Value = 5585.6665
$ws.Range("J61").Value = 3502
$ws.Range("K61").Value = 5585.6665
$ws.Range("L61").Value = 3502
$ws.Range("M61").Value = -5383.6665
$ws.Range("N61").Value = -3906
$ws.Range("H93").Value = 998.3333
$ws.Range("I93").Value = 998.3333
$ws.Range("K93").Value = 998.3333
$ws.Range("M93").Value = 249.6667
$ws.Range("H110").Value = 49994
$ws.Range("J110").Value = 49994
$ws.Range("L110").Value = 49994
$ws.Range("N110").Value = -58174
$ws.Range("H113").Value = 5206.8184
$ws.Range("I113").Value = 5585.6665
$ws.Range("J113").Value = 3502
$ws.Range("K113").Value = 5585.6665
$ws.Range("L113").Value = 3502
$ws.Range("M113").Value = -3415.6665
$ws.Range("N113").Value = -7842
$ws.Range("H132").Value = 3834.2083
$ws.Range("I132").Value = 2884.8667
$ws.Range("J132").Value = 5416.4443
$ws.Range("K132").Value = 8654.6001
$ws.Range("L132").Value = 16249.3329
$ws.Range("M132").Value = -6124.6001
$ws.Range("N132").Value = -21309.3329

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 608.75
$ws.Range("I107").Value = 420.66666
$ws.Range("K107").Value = 1261.99998
$ws.Range("M107").Value = 658.0000199999999
$ws.Range("H122").Value = 10007629
$ws.Range("I122").Value = 10007629
$ws.Range("K122").Value = 30022887
$ws.Range("M122").Value = -30020437
